# Solved 74. Search a 2D Matrix
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new row of data (row 10) for the newly solved problem, following
# the same pattern as the existing "Binary Search" category rows.
$ws.Range("A10").Value = "Binary Search"
$ws.Range("B10").Value = "Search a 2D Matrix"
$ws.Range("C10").Value = "discrete binary search, 2D arrays;"

# Match the alternating highlight formatting used on column B for this
# category block (same formatting as B7/B9) by copying formats only.
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to where the user ended up after finishing
# the edit.
$ws.Range("B12").Select()
